# Generate Report for Handback
# Updates the existing handed-back file (484f6901... -> 2400429e...) and
# appends a new handed-back file (baf86ed3...) as a new row on each of the
# three sheets: Overview, zh-cn, de-de.

$wb = $excel.ActiveWorkbook

$oldGuid = "484f6901-4581-47a0-b030-22d7ae1b9968"
$updGuid = "2400429e-56e2-42a6-b66f-797c45f62aee"
$newGuid = "baf86ed3-4bd3-4514-bb26-9d8583d02b36"

$updHash = "a4da4dcfd12a455ba59d6f5aa5ec3ee881bb413e"
$newHash = "27ec0174b1f5dbc42f6ef2da2dc22807fd0aa9d5"

$dateFmt = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

# -- update row 2 (existing handed-back file) --
$ws1.Range("A2").Value = "$updGuid.md"
$ws1.Range("B2").Hyperlinks.Delete()
$ws1.Hyperlinks.Add(
    $ws1.Range("B2"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3e479601223c5b8a61140f38f253c5b91a12da47/e2e/$updGuid.md",
    [Type]::Missing,
    [Type]::Missing,
    "e2e\$updGuid.md"
) | Out-Null
$ws1.Range("G2").Value = "2016-08-31 13:18:07"
$ws1.Range("G2").NumberFormat = $dateFmt

# -- append row 3 (new handed-back file) via the table so the table/autofilter grows --
$lo1 = $ws1.ListObjects.Item(1)
$lo1.ListRows.Add() | Out-Null

$ws1.Range("A3").Value = "$newGuid.md"
$ws1.Hyperlinks.Add(
    $ws1.Range("B3"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3e479601223c5b8a61140f38f253c5b91a12da47/e2e/$newGuid.md",
    [Type]::Missing,
    [Type]::Missing,
    "e2e\$newGuid.md"
) | Out-Null
$ws1.Range("C3").Value = ".md"
$ws1.Range("E3").Value = "Handed back: in sync with en-US"
$ws1.Range("F3").Value = "Handed back: in sync with en-US"
$ws1.Range("G3").Value = "2016-08-31 13:18:07"
$ws1.Range("G3").NumberFormat = $dateFmt

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

# -- update row 2 (existing handed-back file) --
$ws2.Range("A2").Hyperlinks.Delete()
$ws2.Hyperlinks.Add(
    $ws2.Range("A2"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3e479601223c5b8a61140f38f253c5b91a12da47/e2e/$updGuid.md",
    [Type]::Missing,
    [Type]::Missing,
    "$updGuid.md"
) | Out-Null
$ws2.Range("G2").Value = "$updGuid.$updHash.zh-cn.xlf"
$ws2.Range("H2").Value = "2016-08-31 13:17:56"
$ws2.Range("H2").NumberFormat = $dateFmt
$ws2.Range("I2").Hyperlinks.Delete()
$ws2.Hyperlinks.Add(
    $ws2.Range("I2"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/b960cd8600fdacaaa97cf30ea62fcd51f4140b70/e2e/$updGuid.md",
    [Type]::Missing,
    [Type]::Missing,
    "$updGuid.md"
) | Out-Null
$ws2.Range("J2").Value = "$updGuid.$updHash.zh-cn.xlf"
$ws2.Range("K2").Value = "2016-08-31 13:18:32"
$ws2.Range("K2").NumberFormat = $dateFmt

# -- append row 3 (new handed-back file) --
$lo2 = $ws2.ListObjects.Item(1)
$lo2.ListRows.Add() | Out-Null

$ws2.Hyperlinks.Add(
    $ws2.Range("A3"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3e479601223c5b8a61140f38f253c5b91a12da47/e2e/$newGuid.md",
    [Type]::Missing,
    [Type]::Missing,
    "$newGuid.md"
) | Out-Null
$ws2.Range("B3").Value = ".md"
$ws2.Range("C3").Value = "Handed back: in sync with en-US"
$ws2.Range("D3").Value = "e2e"
$ws2.Range("E3").Value = "ht"
$ws2.Range("F3").Value = "True"
$ws2.Range("G3").Value = "$newGuid.$newHash.zh-cn.xlf"
$ws2.Range("H3").Value = "2016-08-31 13:17:56"
$ws2.Range("H3").NumberFormat = $dateFmt
$ws2.Hyperlinks.Add(
    $ws2.Range("I3"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/b960cd8600fdacaaa97cf30ea62fcd51f4140b70/e2e/$newGuid.md",
    [Type]::Missing,
    [Type]::Missing,
    "$newGuid.md"
) | Out-Null
$ws2.Range("J3").Value = "$newGuid.$newHash.zh-cn.xlf"
$ws2.Range("K3").Value = "2016-08-31 13:18:32"
$ws2.Range("K3").NumberFormat = $dateFmt
$ws2.Range("L3").Value = ""
$ws2.Range("M3").Value = "True"
$ws2.Range("N3").Value = ""
$ws2.Range("O3").Value = "False"
$ws2.Range("P3").Value = ""

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

# -- update row 2 (existing handed-back file) --
$ws3.Range("A2").Hyperlinks.Delete()
$ws3.Hyperlinks.Add(
    $ws3.Range("A2"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3e479601223c5b8a61140f38f253c5b91a12da47/e2e/$updGuid.md",
    [Type]::Missing,
    [Type]::Missing,
    "$updGuid.md"
) | Out-Null
$ws3.Range("G2").Value = "$updGuid.$updHash.de-de.xlf"
$ws3.Range("H2").Value = "2016-08-31 13:18:07"
$ws3.Range("H2").NumberFormat = $dateFmt
$ws3.Range("I2").Hyperlinks.Delete()
$ws3.Hyperlinks.Add(
    $ws3.Range("I2"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/29975b69a94a7990a4d8d6e12aa40945819e6ded/e2e/$updGuid.md",
    [Type]::Missing,
    [Type]::Missing,
    "$updGuid.md"
) | Out-Null
$ws3.Range("J2").Value = "$updGuid.$updHash.de-de.xlf"
$ws3.Range("K2").Value = "2016-08-31 13:18:40"
$ws3.Range("K2").NumberFormat = $dateFmt

# -- append row 3 (new handed-back file) --
$lo3 = $ws3.ListObjects.Item(1)
$lo3.ListRows.Add() | Out-Null

$ws3.Hyperlinks.Add(
    $ws3.Range("A3"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3e479601223c5b8a61140f38f253c5b91a12da47/e2e/$newGuid.md",
    [Type]::Missing,
    [Type]::Missing,
    "$newGuid.md"
) | Out-Null
$ws3.Range("B3").Value = ".md"
$ws3.Range("C3").Value = "Handed back: in sync with en-US"
$ws3.Range("D3").Value = "e2e"
$ws3.Range("E3").Value = "ht"
$ws3.Range("F3").Value = "True"
$ws3.Range("G3").Value = "$newGuid.$newHash.de-de.xlf"
$ws3.Range("H3").Value = "2016-08-31 13:18:07"
$ws3.Range("H3").NumberFormat = $dateFmt
$ws3.Hyperlinks.Add(
    $ws3.Range("I3"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/29975b69a94a7990a4d8d6e12aa40945819e6ded/e2e/$newGuid.md",
    [Type]::Missing,
    [Type]::Missing,
    "$newGuid.md"
) | Out-Null
$ws3.Range("J3").Value = "$newGuid.$newHash.de-de.xlf"
$ws3.Range("K3").Value = "2016-08-31 13:18:40"
$ws3.Range("K3").NumberFormat = $dateFmt
$ws3.Range("L3").Value = ""
$ws3.Range("M3").Value = "True"
$ws3.Range("N3").Value = ""
$ws3.Range("O3").Value = "False"
$ws3.Range("P3").Value = ""

Write-Host "Handback report updated: added $newGuid.md, refreshed $updGuid.md"
